$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    'D2' = '42.244.78'
    'E2' = '  -0.59%  '
    'D3' = '2.231.10'
    'E3' = '  -0.63%  '
    'E4' = '  -0.28%  '
    'D5' = '243.09'
    'E5' = '  -0.82%  '
    'E6' = '  -0.41%  '
    'D7' = '74.33'
    'E7' = '  -1.73%  '
    'E8' = '  +0.02%  '
    'D9' = '0.605'
    'E9' = '  -2.73%  '
    'D10' = '42.76'
    'E10' = '  -2.26%  '
    'D11' = '0.0961'
    'E11' = '  +1.21%  '
    'E12' = '  -3.15%  '
    'D13' = '0.104'
    'E13' = '  +0.26%  '
    'D14' = '2.566.24'
    'E14' = '  -0.49%  '
    'D15' = '14.32'
    'E15' = '  -1.76%  '
    'D16' = '0.838'
    'E16' = '  -2.54%  '
    'D17' = '2.249.31'
    'E17' = '  -0.71%  '
    'D18' = '42.059.56'
    'E18' = '  -0.63%  '
    'E19' = '  +3.90%  '
    'D20' = '6.23'
    'E20' = '  +0.41%  '
    'D21' = '72.86'
    'E21' = '  +1.16%  '
    'D22' = '11.15'
    'E22' = '  +1.60%  '
    'D23' = '230.76'
    'E23' = '  -0.47%  '
    'E24' = '  -5.75%  '
    'E25' = '  +0.01%  '
    'E26' = '  -2.79%  '
    'E27' = '  -0.18%  '
    'E28' = '  -1.07%  '
    'E29' = '  -2.60%  '
    'D30' = '167.07'
    'E30' = '  -0.02%  '
    'D31' = '20.62'
    'E31' = '  -0.47%  '
    'D32' = '5.65'
    'E32' = '  -7.75%  '
    'E33' = '  -1.51%  '
    'D34' = '30.06'
    'E34' = '  -1.58%  '
    'E35' = '  -0.52%  '
    'E36' = '  -7.25%  '
    'E37' = '  -6.87%  '
    'E38' = '  -3.53%  '
    'D39' = '13.22'
    'E39' = '  -3.41%  '
    'E40' = '  -2.05%  '
    'B41' = 'THORChain'
    'C41' = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
    'D41' = '5.71'
    'E41' = '  -0.39%  '
    'B42' = 'MultiversX'
    'C42' = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
    'D42' = '65.16'
    'E42' = '  +2.16%  '
    'E43' = '  -1.36%  '
    'D44' = '8.74'
    'E44' = '  -1.23%  '
    'D45' = '104.68'
    'E45' = '  -1.89%  '
    'E46' = '  -2.09%  '
    'B47' = 'ARBITRUM'
    'C47' = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
    'D47' = '1.12'
    'E47' = '  -1.73%  '
    'B48' = 'NEARProtocol'
    'C48' = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
    'D48' = '2.36'
    'E48' = '  -2.45%  '
    'E49' = '  -0.86%  '
    'E50' = '  -1.40%  '
    'D51' = '2.435.69'
    'E51' = '  -0.80%  '
}

foreach ($cell in $updates.Keys) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $updates[$cell]
    $rng.Style = "Normal"
}
